# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.255.91"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.841.63"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'0.9990"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'241.51"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'0.6703"
$ws.Range("E6").Value = "  -2.41%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.07433"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "'0.2963"
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("D10").Value = "'22.86"
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'0.07719"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").Value = "'5.034"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "'0.6794"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").Value = "'86.39"
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").Value = "1.614.29"
$ws.Range("E15").Value = "  -12.23%  "
$ws.Range("D16").Value = "'6.198"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "'0.000008274"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "28.705.17"
$ws.Range("E18").Value = "  -2.09%  "
$ws.Range("D19").Value = "'229.09"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").Value = "'12.56"
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").Value = "'0.9995"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'7.258"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").Value = "'0.9998"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'160.40"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'8.719"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("D26").Value = "'0.1414"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").Value = "'18.04"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'1.513"
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").Value = "'4.205"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "'4.088"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Value = "'1.187"
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'0.05369"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("B33").Value = "LidoDAOToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D33").Value = "'1.878"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7551"
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "'1.141"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "1.330.86"
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "'0.01807"
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "'0.9242"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").Value = "'5.983"
$ws.Range("E41").Value = "  +5.88%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "'103.28"
$ws.Range("E43").Value = "  -2.26%  "
$ws.Range("D44").Value = "'0.00000000126"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "'0.07771"
$ws.Range("E45").Value = "  +11.40%  "
$ws.Range("D46").Value = "'0.5161"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'1.770"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").Value = "'64.03"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.291"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "1.868.95"
$ws.Range("E50").Value = "  -6.05%  "
$ws.Range("D51").Value = "'0.05924"
$ws.Range("E51").Value = "  -0.01%  "
